$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.345.88"
$ws.Range("E2").Value = "  +3.74%  "
$ws.Range("D3").Value = "2.620.23"
$ws.Range("E3").Value = "  +3.11%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "603.60"
$ws.Range("E5").Value = "  +0.35%  "
$ws.Range("D6").Value = "178.56"
$ws.Range("E6").Value = "  +0.55%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "0.526"
$ws.Range("E8").Value = "  +1.28%  "
$ws.Range("E9").Value = "  +8.83%  "
$ws.Range("D10").Value = "2.618.83"
$ws.Range("E10").Value = "  +3.09%  "
$ws.Range("E11").Value = "  +0.95%  "
$ws.Range("D12").Value = "0.354"
$ws.Range("E12").Value = "  +2.45%  "
$ws.Range("E13").Value = "  +0.15%  "
$ws.Range("D14").Value = "3.102.63"
$ws.Range("E14").Value = "  +4.03%  "
$ws.Range("E15").Value = "  +2.90%  "
$ws.Range("D16").Value = "72.288.86"
$ws.Range("E16").Value = "  +3.77%  "
$ws.Range("E17").Value = "  +1.93%  "
$ws.Range("D18").Value = "2.625.05"
$ws.Range("E18").Value = "  +4.27%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "380.44"
$ws.Range("E19").Value = "  +4.23%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "11.60"
$ws.Range("E20").Value = "  +4.34%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "7.86"
$ws.Range("E21").Value = "  +2.10%  "
$ws.Range("D22").Value = "4.19"
$ws.Range("E22").Value = "  +1.97%  "
$ws.Range("D23").Value = "2.04"
$ws.Range("E23").Value = "  +18.96%  "
$ws.Range("D24").Value = "73.37"
$ws.Range("E24").Value = "  +3.86%  "
$ws.Range("D26").Value = "4.39"
$ws.Range("E26").Value = "  +2.70%  "
$ws.Range("D27").Value = "9.98"
$ws.Range("E27").Value = "  +8.49%  "
$ws.Range("D28").Value = "2.754.64"
$ws.Range("E28").Value = "  +3.29%  "
$ws.Range("E29").Value = "  +0.26%  "
$ws.Range("D30").Value = "0.0₃0953"
$ws.Range("E30").Value = "  +2.77%  "
$ws.Range("E31").Value = "  +4.05%  "
$ws.Range("D32").Value = "518.23"
$ws.Range("E32").Value = "  +0.84%  "
$ws.Range("E33").Value = "  +5.22%  "
$ws.Range("E34").Value = "  +1.89%  "
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").Value = "165.47"
$ws.Range("E36").Value = "  +2.44%  "
$ws.Range("D37").Value = "19.32"
$ws.Range("E37").Value = "  +2.66%  "
$ws.Range("E38").Value = "  +5.33%  "
$ws.Range("D39").Value = "19.05"
$ws.Range("E39").Value = "  +0.72%  "
$ws.Range("E40").Value = "  -7.76%  "
$ws.Range("E41").Value = "  +5.63%  "
$ws.Range("D42").Value = "5.08"
$ws.Range("E42").Value = "  +4.08%  "
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").Value = "2.57"
$ws.Range("E44").Value = "  +4.99%  "
$ws.Range("D45").Value = "0.332"
$ws.Range("E45").Value = "  +2.57%  "
$ws.Range("D46").Value = "39.50"
$ws.Range("E46").Value = "  +1.92%  "
$ws.Range("D47").Value = "149.57"
$ws.Range("E47").Value = "  -1.18%  "
$ws.Range("E48").Value = "  +2.17%  "
$ws.Range("E49").Value = "  +4.30%  "
$ws.Range("E50").Value = "  +5.80%  "
$ws.Range("E51").Value = "  +3.22%  "
